# Update xls to new $ref API.
# The workbook stores cross-references as plain strings of the form
#   link:Namespace.Object.Field
# The new format uses:
#   $ref:Namespace:Object.Field
# (i.e. "link:" -> "$ref:" and the first "." becomes ":")
#
# This script updates every cell that contains such a reference string,
# on both worksheets, and restores the view/selection state recorded in
# the edited workbook (DataBlocks sheet active/selected, new selections).

$wb = $excel.ActiveWorkbook

$tests = $wb.Worksheets.Item("Tests")
$blocks = $wb.Worksheets.Item("DataBlocks")

# --- Tests sheet (sheet1) ---
$tests.Range("C5").Value = '$ref:Tests:Uncommon.reftest'
$tests.Range("C6").Value = '$ref:Tests:Uncommon.reftestGen'
$tests.Range("C9").Value = '$ref:DataBlocks:AnotherObject.cyclicRef'
$tests.Range("C10").Value = '$ref:DataBlocks:NewObject'
$tests.Range("C14").Value = '$ref:DataBlocks:NewObject.VALUE'
$tests.Range("C22").Value = '$ref:DataBlocks:NewObject'

# --- DataBlocks sheet (sheet2) ---
$blocks.Range("C3").Value = '$ref:DataBlocks:AnotherObject.anotherValue'
$blocks.Range("C7").Value = '$ref:Tests:Common.cyclic'
$blocks.Range("C8").Value = '$ref:Tests:Common.gendata'

# --- Restore view/selection state recorded in the edited workbook ---
[void]$tests.Range("C5").Select()
[void]$blocks.Activate()
[void]$blocks.Range("C4").Select()
